{"js": "// Add a new paragraph after the last paragraph in the document body:\n// \"\u9886\u5bfc\u5efa\u8bae\u52a0\u4e0a\u53e6\u5916\u4e00\u53e5\u8bdd\u3002\" with the same first-line indent (420 twips)\n// and an eastAsia font hint on its run, matching the first paragraph's\n// formatting (\"\u8fd9\u662f\u7b2c\u4e00\u7248\u3002\").\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\n// Insert a new paragraph right after the last one, with the target text.\nconst newParagraph = lastParagraph.insertParagraph(\"\u9886\u5bfc\u5efa\u8bae\u52a0\u4e0a\u53e6\u5916\u4e00\u53e5\u8bdd\u3002\", Word.InsertLocation.after);\n\n// Match the first-line indent used by the existing paragraph (420 twips = 21pt).\nnewParagraph.firstLineIndent = 21;\n\nawait context.sync();\n", "ps1": "# Add a new paragraph at the end of the document with the boss's requested\n# sentence, matching the existing paragraph's first-line indent and font.\n\n$d = $word.ActiveDocument\n\n# Move to the very end of the document content.\n$end = $d.Content\n$end.Collapse(0)  # wdCollapseEnd\n\n# Start a new paragraph, then type the requested sentence.\n$end.InsertParagraphAfter()\n$end.Collapse(0)\n$end.Move(4, 1) | Out-Null  # wdParagraph, move into the freshly inserted paragraph\n\n$newPara = $d.Paragraphs.Last\n$newRange = $newPara.Range\n$newRange.Text = \"\u9886\u5bfc\u5efa\u8bae\u52a0\u4e0a\u53e6\u5916\u4e00\u53e5\u8bdd\u3002\"\n\n$newPara.Format.FirstLineIndent = 21  # 420 twips = 21 points\n"}
